# Refresh the cryptos price/volume snapshot (Sheet1, columns D "Price" and
# E "Volume(1h)") to the latest scraped values.
#
# Several new "Price" strings (e.g. "1.00", "27.99") are digit/period-only
# and would otherwise be auto-converted from Text to a Number by Excel on
# assignment; the source data models these as plain text, so those specific
# cells are written with a leading apostrophe (forces Text) and then have
# their Style reset to "Normal" so no stray quote-prefix style is left
# behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.376.71'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '2.646.66'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'597.69"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = "'158.89"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.89%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('E9').Value = '  -1.06%  '
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').Value = "'27.99"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('E15').Value = '  -2.80%  '
$ws.Range('D16').Value = '68.362.59'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '2.636.51'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = "'11.37"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').Value = "'359.67"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.20%  '
$ws.Range('E20').Value = '  -1.19%  '
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('D22').Value = "'4.76"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.61%  '
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').Value = "'74.53"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = "'1.00"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  -0.75%  '
$ws.Range('E28').Value = '  -2.62%  '
$ws.Range('D29').Value = "'0.998"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('D30').Value = "'562.17"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('E31').Value = '  -1.52%  '
$ws.Range('D32').Value = "'1.39"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.05%  '
$ws.Range('E33').Value = '  +1.11%  '
$ws.Range('D34').Value = "'1.64"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.94%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -1.75%  '
$ws.Range('D37').Value = "'160.14"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('D38').Value = "'19.67"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('D40').Value = "'1.86"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').Value = "'5.33"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.49%  '
$ws.Range('D42').Value = "'2.62"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.38%  '
$ws.Range('D43').Value = '0.0₆0320'
$ws.Range('E43').Value = '  -4.74%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = "'157.41"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('D46').Value = "'3.80"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.58%  '
$ws.Range('D47').Value = "'22.01"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('E48').Value = '  -1.47%  '
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('E50').Value = '  +1.16%  '
$ws.Range('E51').Value = '  -0.22%  '
